# Update figures: split the "raw"/"formatted" pair of sheets into two pairs
# — one for "all WRGs" (the existing data, renamed) and a new one for
# "simple WRGs" (new raw data + a formatted table copied from the original).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "copy raw data here"      -> "raw - all WRGs"
$ws2 = $wb.Worksheets.Item(2)   # "Formatted table"          -> "formatted - all WRGs"

# Rename the original two sheets. Excel keeps cross-sheet formula references
# in sync automatically when a sheet is renamed.
$ws1.Name = "raw - all WRGs"
$ws2.Name = "formatted - all WRGs"

# Duplicate the raw + formatted sheets to build the "simple WRGs" pair,
# keeping them positioned right after the sheets they were copied from.
$ws1.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "raw - simple WRGs"

$ws2.Copy($null, $ws3)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "formatted - simple WRGs"

# --- Populate "raw - simple WRGs" with the new raw data ---------------------
# Columns: A = Algorithm (unchanged), B = Bias_prc, C = MAE_1e5m3, D = R2, E = slope

$ws3.Range("B2").Value = -35.5
$ws3.Range("C2").Value = 0.64236578016942403
$ws3.Range("D2").Value = 0.27127328130709299
$ws3.Range("E2").Value = 0.73073188677714296

$ws3.Range("B3").Value = 3.9
$ws3.Range("C3").Value = 0.64100551487593405
$ws3.Range("D3").Value = 0.27375581649763697
$ws3.Range("E3").Value = 0.99439206244612299

$ws3.Range("B4").Value = -20.9
$ws3.Range("C4").Value = 0.61820686485772802
$ws3.Range("D4").Value = 0.27981023355653301
$ws3.Range("E4").Value = 0.84696922688904497

$ws3.Range("B5").Value = -45.9
$ws3.Range("C5").Value = 0.77320462929236899
$ws3.Range("D5").Value = 0.209262556050543
$ws3.Range("E5").Value = 0.67818615166065499

$ws3.Range("B6").Value = -31
$ws3.Range("C6").Value = 0.66123445663315195
$ws3.Range("D6").Value = 0.25439071056750601
$ws3.Range("E6").Value = 0.76193852997043998

$ws3.Range("B7").Value = 33.9
$ws3.Range("C7").Value = 0.63872081186159302
$ws3.Range("D7").Value = 0.29571647620309
$ws3.Range("E7").Value = 0.95811764316371195

$ws3.Range("B8").Value = -23.6
$ws3.Range("C8").Value = 0.61178393848644297
$ws3.Range("D8").Value = 0.29325253206296797
$ws3.Range("E8").Value = 0.84286612540657502

# --- Point "formatted - simple WRGs" formulas at the new raw sheet ---------
# (keeps the same column-reshuffle pattern as the "all WRGs" formatted sheet:
#  B<-C, C<-B, D<-E, E<-D)

for ($r = 2; $r -le 8; $r++) {
    $ws4.Range("B$r").Formula = "='raw - simple WRGs'!C$r"
    $ws4.Range("C$r").Formula = "='raw - simple WRGs'!B$r"
    $ws4.Range("D$r").Formula = "='raw - simple WRGs'!E$r"
    $ws4.Range("E$r").Formula = "='raw - simple WRGs'!D$r"
}

# --- View state: the new "formatted - simple WRGs" sheet becomes the active
# tab, and the other sheets' selections are no longer the one shown when the
# file opens.
$ws1.Range("G24").Select()
$ws2.Range("D14").Select()
$ws3.Cells.Select()
$ws4.Range("G24").Select()

$wb.Windows.Item(1).DisplayedTab = $ws4
$excel.ActiveWindow.ActiveSheet = $ws4
